$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 66, shifting existing rows 66:166 down to 67:167
$ws.Rows.Item(66).Insert()

# Populate the new row 66 with the new data record
$ws.Range("A66").Value = 9
$ws.Range("B66").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C66").Value = "Metropolitana"
$ws.Range("D66").Value = 44973
$ws.Range("E66").Value = 13
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100103
$ws.Range("H66").Value = "Frutos de hueso (carozo)"
$ws.Range("I66").Value = 100103002
$ws.Range("J66").Value = "Ciruela"
$ws.Range("K66").Value = "Friar"
$ws.Range("L66").Value = "Primera"
$ws.Range("M66").Value = 5
$ws.Range("N66").Value = 190000
$ws.Range("O66").Value = 190000
$ws.Range("P66").Value = 190000
$ws.Range("Q66").Value = "$/bins (450 kilos)"
$ws.Range("R66").Value = "Región Metropolitana"
$ws.Range("S66").Value = 422
$ws.Range("T66").Value = 450

# Ensure date cell keeps proper numeric formatting matching the rest of column D
$ws.Range("D66").NumberFormat = $ws.Range("D67").NumberFormat
